# Insert a new weekly price record at row 476, shifting the existing
# rows 476-531 down to 477-532 (dimension grows from R531 to R532).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("476:476").Insert()

$ws.Cells.Item(476, 1).Value  = 4
$ws.Cells.Item(476, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(476, 3).Value  = "Los Lagos"
$ws.Cells.Item(476, 4).Value  = 45194
$ws.Cells.Item(476, 5).Value  = 10
$ws.Cells.Item(476, 6).Value  = 100112040
$ws.Cells.Item(476, 7).Value  = "Cilantro"
$ws.Cells.Item(476, 8).Value  = "Sin especificar"
$ws.Cells.Item(476, 9).Value  = "Primera"
$ws.Cells.Item(476, 10).Value = 80
$ws.Cells.Item(476, 11).Value = 12000
$ws.Cells.Item(476, 12).Value = 12000
$ws.Cells.Item(476, 13).Value = 12000
$ws.Cells.Item(476, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(476, 15).Value = "Región Metropolitana"
$ws.Cells.Item(476, 16).Value = 333
$ws.Cells.Item(476, 17).Value = 36
$ws.Cells.Item(476, 18).Value = "Hortaliza"

# Keep the date cell formatted like the rest of column D.
$ws.Cells.Item(476, 4).NumberFormat = $ws.Cells.Item(477, 4).NumberFormat
